$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parts Purchased")
$ws.Rows.Item(26).Insert()
